$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46, shifting existing rows 46:76 down to 47:77.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new reading
# (same dimension/category columns as its neighbours, new date/variety/price data).
$ws.Range("A46").Value = 7
$ws.Range("B46").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C46").Value = "Ñuble"
$ws.Range("D46").Value = 44960
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100103
$ws.Range("H46").Value = "Frutos de hueso (carozo)"
$ws.Range("I46").Value = 100103002
$ws.Range("J46").Value = "Ciruela"
$ws.Range("K46").Value = "Fortuna"
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 60
$ws.Range("N46").Value = 11000
$ws.Range("O46").Value = 12000
$ws.Range("P46").Value = 11500
$ws.Range("Q46").Value = "$/bandeja 18 kilos granel"
$ws.Range("R46").Value = "Región de O'Higgins"
$ws.Range("S46").Value = 639
$ws.Range("T46").Value = 18
